# "Añadido visor control. Más fácil depurar ahora."
# Reworks the "Registros pipeline" block (adds a K/L "raw value -> relative"
# pair fed straight off the new "Banco de registros" lookup table in I6:I9,
# renames the old M/N columns to O/P), relocates the second block
# ("De relativos a absolutos") a few rows down, and adds three new
# mini "viewer/control" calculators below it: Banco de registros (25-30),
# NPCs (32-34) and Control (37-45).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 5: new header label next to "Registros pipeline"
# ---------------------------------------------------------------------
$ws.Range("I5").Value = "Derecho Buffer Instr"

# ---------------------------------------------------------------------
# Rows 6-9: "Registros pipeline" table rework
# ---------------------------------------------------------------------
# old E6 ("Objetivo 1" label) is dropped - no longer used here
$ws.Range("E6").ClearContents()

# K/L is a straight copy of what I/J used to contain (fixed reference
# values), I/J now pulls live off the new "Banco de registros" block
$ws.Range("K6").Value = 242
$ws.Range("K7").Value = 260
$ws.Range("K8").Value = 330
$ws.Range("K9").Value = 345

$ws.Range("L6").Formula = '=(K6-$B$6)/$B$8'
$ws.Range("L7").Formula = '=(K7-$B$6)/$B$8'
$ws.Range("L8").Formula = '=(K8-$B$6)/$B$8'
$ws.Range("L9").Formula = '=(K9-$B$6)/$B$8'

$ws.Range("I6").Formula = '=G25'
$ws.Range("I7").Formula = '=G26'
$ws.Range("I8").Formula = '=G27'
$ws.Range("I9").Formula = '=G28'

# J6:J9 keep their original formula - value only changes because I6:I9 did
$ws.Range("J6").Formula = '=(I6-$B$6)/$B$8'
$ws.Range("J7").Formula = '=(I7-$B$6)/$B$8'
$ws.Range("J8").Formula = '=(I8-$B$6)/$B$8'
$ws.Range("J9").Formula = '=(I9-$B$6)/$B$8'

# M/N (rows 6-8) move over to O/P
$ws.Range("O6").Value = $ws.Range("M6").Value2
$ws.Range("O7").Value = $ws.Range("M7").Value2
$ws.Range("O8").Value = $ws.Range("M8").Value2
$ws.Range("M6:N8").ClearContents()

$ws.Range("P6").Formula = '=(O6-$B$6)/$B$8'
$ws.Range("P7").Formula = '=(O7-$B$6)/$B$8'
$ws.Range("P8").Formula = '=(O8-$B$6)/$B$8'

Write-Host "Rows 5-9 reworked"

# ---------------------------------------------------------------------
# "De relativos a absolutos" block: was rows 10/12-14, now rows 16/18-20
# ---------------------------------------------------------------------
$ws.Range("A16").Value = $ws.Range("A10").Value2
$ws.Range("A10").ClearContents()

$ws.Range("A18").Value = $ws.Range("A12").Value2
$ws.Range("B18").Value = $ws.Range("B12").Value2
$ws.Range("E18").Value = $ws.Range("E12").Value2
$ws.Range("F18").Value = $ws.Range("F12").Value2
$ws.Range("G18").Formula = '=$B$18+F18*$B$20'

$ws.Range("A19").Value = $ws.Range("A13").Value2
$ws.Range("B19").Value = $ws.Range("B13").Value2
$ws.Range("F19").Value = $ws.Range("F13").Value2
$ws.Range("G19").Formula = '=$B$18+F19*$B$20'

$ws.Range("A20").Value = $ws.Range("A14").Value2
$ws.Range("B20").Formula = '=B19-B18'

$ws.Range("A12:G14").ClearContents()

Write-Host "Moved rows 10-14 block to 16-20"

# ---------------------------------------------------------------------
# New "Banco de registros" mini-table (rows 25-30), shared formula in G
# ---------------------------------------------------------------------
$ws.Range("A25").Value = "Banco de registros"
$ws.Range("B25").Value = 100
$ws.Range("F25").Value = 0.2
$ws.Range("G25:G30").Formula = '=$B$18+F25*$B$20'

$ws.Range("B26").Value = 220
$ws.Range("F26").Value = 0.4

$ws.Range("B27").Formula = '=B26-B25'
$ws.Range("F27").Value = 0.6

$ws.Range("F28").Value = 0.8
$ws.Range("F29").Value = 0.25
$ws.Range("F30").Value = 0.65

Write-Host "Added Banco de registros block"

# ---------------------------------------------------------------------
# New "NPCs" mini calculator (rows 32-34)
# ---------------------------------------------------------------------
$ws.Range("A32").Value = "NPCs"
$ws.Range("B32").Value = 120
$ws.Range("F32").Value = 150
$ws.Range("G32").Formula = '=(F32-$B$32)/$B$34'
$ws.Range("B33").Value = 162
$ws.Range("B34").Value = 42

Write-Host "Added NPCs block"

# ---------------------------------------------------------------------
# New "Control" viewer (rows 37-45), shared formula in G
# ---------------------------------------------------------------------
$ws.Range("A37").Value = "Control"
$ws.Range("B37").Value = 275
$ws.Range("F37").Value = 480
$ws.Range("G37").Formula = '=(F37-$B$37)/$B$39'

$ws.Range("F38").Value = 560
$ws.Range("G38:G45").Formula = '=(F38-$B$37)/$B$39'
$ws.Range("B39").Value = 1070
$ws.Range("F39").Value = 615
$ws.Range("F40").Value = 660
$ws.Range("F41").Value = 870
$ws.Range("F42").Value = 950
$ws.Range("F43").Value = 1050
$ws.Range("F44").Value = 1140
$ws.Range("F45").Value = 1316.66

Write-Host "Added Control block"

# ---------------------------------------------------------------------
# View: scroll to the new Control block and select its result column
# ---------------------------------------------------------------------
$ws.Range("G37:G45").Select()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1

Write-Host "Done"
